# Update "想去人数" (want-to-go count) values in column F across sheets,
# reflecting the latest scrape output.

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 13914
$ws1.Range("F3").Value = 328
$ws1.Range("F4").Value = 671
$ws1.Range("F6").Value = 511
$ws1.Range("F7").Value = 1440

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 10

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 13914
$ws4.Range("F3").Value = 328
$ws4.Range("F4").Value = 671
$ws4.Range("F6").Value = 10
$ws4.Range("F8").Value = 511
$ws4.Range("F9").Value = 1440
